$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "Lost in Space is an exploration, puzzle solving game, which throws players into the seat of a lost astronaut trying to find his way back home or make a difficult choice to save the people living on this planet and never see his loved ones again.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Architecture Simulator is a simulation which focuses on giving the player the ability to design their own house and walk around in. They can select between different models to interact with and do as they see fit.",
    2
)
